# Add 32 more "Element #" rows (65-96) to the Dongle impedance test sheet,
# pushing the trailing spacer row and the "Notes" footer row further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74 (Element #64) is the last data row before the spacer (row 75) and
# the Notes footer (row 76). Insert 32 blank rows right before the spacer
# row so it - and the footer after it - shift from 75/76 down to 107/108.
$ws.Range("A75:A106").EntireRow.Insert()

# Fill the 32 newly inserted rows (now 75-106) the same way the existing
# data rows (11-74) look: copy the formatting of the last data row (74)
# into each new row, then set the incrementing "Element #" value.
for ($i = 0; $i -lt 32; $i++) {
    $destRow = 75 + $i
    $ws.Range("C74:J74").Copy($ws.Range("C" + $destRow + ":J" + $destRow))
    $ws.Range("C" + $destRow).Value = 65 + $i
}

# The first of the new rows keeps an explicit (default) row height, matching
# the source workbook.
$ws.Rows.Item(75).RowHeight = 14.4

# Reflect the final selection left behind by the edit.
$ws.Range("N107").Select()
